# For every cell in column G ("Recorded By") that contains a
# comma-separated list of recorder names/emails, rotate the list left
# by one position - i.e. move the first entry to the end of the list.
#
# Examples:
#   "System, dnasr281@gmail.com"              -> "dnasr281@gmail.com, System"
#   "System, system, backup@backdoor.com"     -> "system, backup@backdoor.com, System"
#   "admin@admin.com, dnasr281@gmail.com"     -> "dnasr281@gmail.com, admin@admin.com"
#
# Cells with a single value (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($null -ne $val -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ",\s*"
        if ($parts.Count -gt 1) {
            $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
            $cell.Value2 = $rotated
        }
    }
}
